$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")
$ws.Activate()

# Update B2 value (15.9 -> 11)
$ws.Range("B2").Value = 11

# Insert a new row at 19, shifting existing rows 19:53 down to 20:54
$ws.Rows("19:19").Insert()

# Populate the newly inserted row 19 with the new parameter
$ws.Range("A19").Value = "prop_mix_ruralpoor_from_ruralpoor"
$ws.Range("B19").Value = 0.6

# Update the selected cell shown when the workbook is opened
$ws.Range("B3").Select()
